$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the stand-alone "2346890 - Eliane Correa Pedrozo" data
#    row (old row 13, B/C only, no label in column A). Deleting the
#    whole row shifts everything below it up by one.
# ------------------------------------------------------------------
$ws.Rows.Item(13).Delete()

# ------------------------------------------------------------------
# 2) Remove the long Bibliografia reference-list row. After step 1
#    it now lives at row 21. Deleting it shifts everything below up.
# ------------------------------------------------------------------
$ws.Rows.Item(21).Delete()

# ------------------------------------------------------------------
# 3) Insert a fresh blank row before the current row 21 (Requisitos:)
#    to make room for the relocated "Bibliografia:" label + value.
# ------------------------------------------------------------------
$ws.Rows.Item(21).Insert()

# ------------------------------------------------------------------
# 4) Fix up cell contents that moved / were replaced.
# ------------------------------------------------------------------

# Row 10 (Objetivos:) - B/C now hold the docente reference instead of
# the long objectives paragraph.
$ws.Cells.Item(10,2).Value2 = "2346890 - Eliane Corrêa Pedrozo"
$ws.Cells.Item(10,3).Value2 = "2346890 - Eliane Corrêa Pedrozo"

# Row 13 (Programa resumido:) - B/C now hold "Semestral" instead of
# the long summary paragraph.
$ws.Cells.Item(13,2).Value2 = "Semestral"
$ws.Cells.Item(13,3).Value2 = "Semestral"

# Row 15 (Programa:) - gains a B/C value of "15/07/2015". Copy the
# formatting from row 13 (which already carries the proper B/C
# styles) so the shared style indices line up instead of minting new
# duplicate styles.
$ws.Cells.Item(13,2).Copy()
$ws.Cells.Item(15,2).PasteSpecial(-4122)
$ws.Cells.Item(13,3).Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4122)
$ws.Cells.Item(15,2).Value2 = "15/07/2015"
$ws.Cells.Item(15,3).Value2 = "15/07/2015"
$ws.Rows.Item(15).RowHeight = 120

# Row 17 (Avaliação:) no longer carries a custom row height.
$ws.Rows.Item(17).AutoFit()

# Row 18 (Método:) gains a B/C value of the docente reference. Copy
# formatting from row 19 (already B/C-styled) first.
$ws.Cells.Item(19,2).Copy()
$ws.Cells.Item(18,2).PasteSpecial(-4122)
$ws.Cells.Item(19,3).Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4122)
$ws.Cells.Item(18,2).Value2 = "2346890 - Eliane Corrêa Pedrozo"
$ws.Cells.Item(18,3).Value2 = "2346890 - Eliane Corrêa Pedrozo"
$ws.Rows.Item(18).RowHeight = 60

# Row 19 (Critério:) now holds the "Seminário em grupo..." text that
# used to sit on the Método: row.
$ws.Cells.Item(19,2).Value2 = "Seminário em grupo sobre um estudo de caso apresentado pelos alunos.`nResolução de um exercício individual após cada aula sobre o tema abordado, com consulta. `nProva escrita."
$ws.Cells.Item(19,3).Value2 = "Seminário em grupo sobre um estudo de caso apresentado pelos alunos.`nResolução de um exercício individual após cada aula sobre o tema abordado, com consulta. `nProva escrita."

# Row 20 (Norma de recuperação:) now holds the "Média Final..." text
# that used to sit on the Critério: row.
$ws.Cells.Item(20,2).Value2 = "Média Final = 0,4 x Nota da Prova + 0,2 x Nota dos exercícios + 0,4 x Nota do Seminário`n`nMédia Final Mínima para Aprovação = 5,0"
$ws.Cells.Item(20,3).Value2 = "Média Final = 0,4 x Nota da Prova + 0,2 x Nota dos exercícios + 0,4 x Nota do Seminário`n`nMédia Final Mínima para Aprovação = 5,0"

# Row 21 (new) becomes the relocated "Bibliografia:" row, holding the
# "Nota Final..." text that used to sit on the Norma de recuperação:
# row. Copy formatting from row 20 first (A/B/C all already styled).
$ws.Cells.Item(20,1).Copy()
$ws.Cells.Item(21,1).PasteSpecial(-4122)
$ws.Cells.Item(20,2).Copy()
$ws.Cells.Item(21,2).PasteSpecial(-4122)
$ws.Cells.Item(20,3).Copy()
$ws.Cells.Item(21,3).PasteSpecial(-4122)

$ws.Cells.Item(21,1).Value2 = "Bibliografia:"
$ws.Cells.Item(21,2).Value2 = "Nota Final = (Prova Escrita + Média final)/2`n`nNota Final Mínima para Aprovação = 5,0"
$ws.Cells.Item(21,3).Value2 = "Nota Final = (Prova Escrita + Média final)/2`n`nNota Final Mínima para Aprovação = 5,0"
$ws.Rows.Item(21).RowHeight = 120

$ws.Range("A1").Select()
